$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "35.317.27"
$ws.Cells.Item(2, 5).Value = "  +0.23%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.900.01"
$ws.Cells.Item(3, 5).Value = "  +1.97%  "
$ws.Cells.Item(4, 5).Value = "  -0.28%  "
$ws.Cells.Item(5, 5).Value = "  +2.55%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.665"
$ws.Cells.Item(6, 5).Value = "  +6.55%  "
$ws.Cells.Item(7, 5).Value = "  -0.24%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "41.19"
$ws.Cells.Item(8, 5).Value = "  -3.03%  "
$ws.Cells.Item(9, 5).Value = "  +5.57%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "52.97"
$ws.Cells.Item(10, 5).Value = "  +12.96%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0718"
$ws.Cells.Item(11, 5).Value = "  +3.18%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0994"
$ws.Cells.Item(12, 5).Value = "  +0.40%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "2.176.94"
$ws.Cells.Item(13, 5).Value = "  +2.06%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "12.09"
$ws.Cells.Item(14, 5).Value = "  +5.17%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.695"
$ws.Cells.Item(15, 5).Value = "  +2.39%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "1.896.75"
$ws.Cells.Item(16, 5).Value = "  +2.08%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "4.82"
$ws.Cells.Item(17, 5).Value = "  +1.88%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "35.325.06"
$ws.Cells.Item(18, 5).Value = "  +0.24%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "72.13"
$ws.Cells.Item(19, 5).Value = "  +3.06%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.0₃0817"
$ws.Cells.Item(20, 5).Value = "  +2.59%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "240.45"
$ws.Cells.Item(21, 5).Value = "  -0.42%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "12.41"
$ws.Cells.Item(22, 5).Value = "  +1.23%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.82"
$ws.Cells.Item(23, 5).Value = "  +1.49%  "
$ws.Cells.Item(24, 5).Value = "  -0.25%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.29"
$ws.Cells.Item(25, 5).Value = "  +1.35%  "
$ws.Cells.Item(26, 5).Value = "  +22.23%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "170.33"
$ws.Cells.Item(27, 5).Value = "  +0.40%  "
$ws.Cells.Item(28, 5).Value = "  +3.14%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "18.30"
$ws.Cells.Item(29, 5).Value = "  +3.09%  "
$ws.Cells.Item(30, 5).Value = "  +1.60%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.12"
$ws.Cells.Item(31, 5).Value = "  +2.13%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.0563"
$ws.Cells.Item(32, 5).Value = "  +0.06%  "
$ws.Cells.Item(33, 5).Value = "  -0.24%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.927"
$ws.Cells.Item(34, 5).Value = "  +13.29%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.07"
$ws.Cells.Item(35, 5).Value = "  +0.69%  "
$ws.Cells.Item(36, 5).Value = "  -5.56%  "
$ws.Cells.Item(37, 5).Value = "  -2.50%  "
$ws.Cells.Item(38, 5).Value = "  +1.48%  "
$ws.Cells.Item(39, 5).Value = "  -1.02%  "
$ws.Cells.Item(40, 5).Value = "  +2.38%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "16.00"
$ws.Cells.Item(41, 5).Value = "  +4.61%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0624"
$ws.Cells.Item(42, 5).Value = "  +3.40%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "89.29"
$ws.Cells.Item(43, 5).Value = "  -1.24%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.335.58"
$ws.Cells.Item(44, 5).Value = "  -0.78%  "
$ws.Cells.Item(45, 5).Value = "  +0.55%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "47.90"
$ws.Cells.Item(46, 5).Value = "  +37.33%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.41"
$ws.Cells.Item(47, 5).Value = "  -0.66%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "6.48"
$ws.Cells.Item(49, 5).Value = "  -1.79%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "11.81"
$ws.Cells.Item(50, 5).Value = "  -5.22%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.086.27"
$ws.Cells.Item(51, 5).Value = "  +1.92%  "
